$wb = $excel.ActiveWorkbook

# Source range used to copy the header-row style (s="1") onto the new sheets'
# row-1 cells, so we reuse the existing style index instead of creating a new one.
$styleSource = $wb.Worksheets.Item("Backbones")

# ----- Backbones74 (was xl/worksheets/sheet76.xml) -----
$count = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($count))
$ws.Name = "Backbones74"

$styleSource.Range("A1:H1").Copy($ws.Range("A1:H1"))

$ws.Cells.Item(1, 1).Value = 8
$ws.Cells.Item(1, 2).Value = 60
$ws.Cells.Item(1, 3).Value = 88
$ws.Cells.Item(1, 4).Value = 97
$ws.Cells.Item(1, 5).Value = 1
$ws.Cells.Item(1, 6).Value = 95
$ws.Cells.Item(1, 7).Value = 5
$ws.Cells.Item(1, 8).Value = 58
$ws.Cells.Item(2, 1).Value = 23
$ws.Cells.Item(2, 2).Value = 94
$ws.Cells.Item(2, 3).Value = 4
$ws.Cells.Item(2, 4).Value = 40
$ws.Cells.Item(2, 5).Value = 56
$ws.Cells.Item(2, 6).Value = 82
$ws.Cells.Item(2, 7).Value = 33
$ws.Cells.Item(3, 1).Value = 66
$ws.Cells.Item(3, 2).Value = 14
$ws.Cells.Item(3, 3).Value = 96
$ws.Cells.Item(3, 4).Value = 35
$ws.Cells.Item(3, 5).Value = 98
$ws.Cells.Item(3, 6).Value = 27
$ws.Cells.Item(3, 7).Value = 49
$ws.Cells.Item(4, 1).Value = 21
$ws.Cells.Item(4, 2).Value = 36
$ws.Cells.Item(4, 3).Value = 48
$ws.Cells.Item(4, 4).Value = 51
$ws.Cells.Item(4, 5).Value = 100
$ws.Cells.Item(4, 6).Value = 84
$ws.Cells.Item(4, 7).Value = 69
$ws.Cells.Item(5, 1).Value = 70
$ws.Cells.Item(5, 2).Value = 24
$ws.Cells.Item(5, 3).Value = 18
$ws.Cells.Item(5, 4).Value = 62
$ws.Cells.Item(5, 5).Value = 83
$ws.Cells.Item(5, 6).Value = 68
$ws.Cells.Item(5, 7).Value = 10
$ws.Cells.Item(6, 1).Value = 20
$ws.Cells.Item(6, 2).Value = 22
$ws.Cells.Item(6, 4).Value = 86
$ws.Cells.Item(6, 5).Value = 28
$ws.Cells.Item(6, 7).Value = 75
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 53
$ws.Cells.Item(7, 4).Value = 59
$ws.Cells.Item(7, 5).Value = 41
$ws.Cells.Item(7, 7).Value = 57
$ws.Cells.Item(8, 1).Value = 52
$ws.Cells.Item(8, 2).Value = 74
$ws.Cells.Item(8, 4).Value = 26
$ws.Cells.Item(8, 5).Value = 77
$ws.Cells.Item(8, 7).Value = 45
$ws.Cells.Item(9, 1).Value = 92
$ws.Cells.Item(9, 2).Value = 87
$ws.Cells.Item(9, 4).Value = 17
$ws.Cells.Item(9, 5).Value = 79
$ws.Cells.Item(10, 1).Value = 12
$ws.Cells.Item(10, 2).Value = 13
$ws.Cells.Item(10, 4).Value = 38
$ws.Cells.Item(10, 5).Value = 78
$ws.Cells.Item(11, 1).Value = 44
$ws.Cells.Item(11, 2).Value = 37
$ws.Cells.Item(11, 4).Value = 3
$ws.Cells.Item(11, 5).Value = 67
$ws.Cells.Item(12, 1).Value = 9
$ws.Cells.Item(12, 2).Value = 16
$ws.Cells.Item(12, 4).Value = 55
$ws.Cells.Item(12, 5).Value = 64
$ws.Cells.Item(13, 1).Value = 31
$ws.Cells.Item(13, 2).Value = 2
$ws.Cells.Item(13, 4).Value = 43
$ws.Cells.Item(13, 5).Value = 29
$ws.Cells.Item(14, 1).Value = 90
$ws.Cells.Item(14, 2).Value = 32
$ws.Cells.Item(14, 4).Value = 47
$ws.Cells.Item(15, 1).Value = 73
$ws.Cells.Item(15, 2).Value = 99
$ws.Cells.Item(15, 4).Value = 89
$ws.Cells.Item(16, 1).Value = 11
$ws.Cells.Item(16, 2).Value = 71
$ws.Cells.Item(16, 4).Value = 85
$ws.Cells.Item(17, 1).Value = 81
$ws.Cells.Item(17, 2).Value = 91
$ws.Cells.Item(17, 4).Value = 63
$ws.Cells.Item(18, 1).Value = 42
$ws.Cells.Item(18, 2).Value = 7
$ws.Cells.Item(18, 4).Value = 93
$ws.Cells.Item(19, 1).Value = 50
$ws.Cells.Item(19, 2).Value = 54
$ws.Cells.Item(19, 4).Value = 19
$ws.Cells.Item(20, 1).Value = 80
$ws.Cells.Item(20, 2).Value = 39
$ws.Cells.Item(20, 4).Value = 34
$ws.Cells.Item(21, 1).Value = 15
$ws.Cells.Item(21, 2).Value = 25
$ws.Cells.Item(22, 1).Value = 61
$ws.Cells.Item(22, 2).Value = 46
$ws.Cells.Item(23, 1).Value = 76
$ws.Cells.Item(24, 1).Value = 30
$ws.Cells.Item(25, 1).Value = 72
$ws.Cells.Item(26, 1).Value = 65

# ----- Backbones75 (was xl/worksheets/sheet77.xml) -----
$count = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($count))
$ws.Name = "Backbones75"

$styleSource.Range("A1:H1").Copy($ws.Range("A1:H1"))

$ws.Cells.Item(1, 1).Value = 8
$ws.Cells.Item(1, 2).Value = 60
$ws.Cells.Item(1, 3).Value = 88
$ws.Cells.Item(1, 4).Value = 97
$ws.Cells.Item(1, 5).Value = 20
$ws.Cells.Item(1, 6).Value = 46
$ws.Cells.Item(1, 7).Value = 50
$ws.Cells.Item(1, 8).Value = 43
$ws.Cells.Item(2, 1).Value = 53
$ws.Cells.Item(2, 2).Value = 30
$ws.Cells.Item(2, 3).Value = 21
$ws.Cells.Item(2, 5).Value = 34
$ws.Cells.Item(2, 6).Value = 89
$ws.Cells.Item(2, 7).Value = 17
$ws.Cells.Item(2, 8).Value = 38
$ws.Cells.Item(3, 1).Value = 75
$ws.Cells.Item(3, 2).Value = 66
$ws.Cells.Item(3, 3).Value = 28
$ws.Cells.Item(3, 5).Value = 57
$ws.Cells.Item(3, 6).Value = 74
$ws.Cells.Item(3, 8).Value = 37
$ws.Cells.Item(4, 1).Value = 48
$ws.Cells.Item(4, 3).Value = 49
$ws.Cells.Item(4, 5).Value = 52
$ws.Cells.Item(4, 6).Value = 27
$ws.Cells.Item(4, 8).Value = 18
$ws.Cells.Item(5, 1).Value = 76
$ws.Cells.Item(5, 3).Value = 86
$ws.Cells.Item(5, 5).Value = 11
$ws.Cells.Item(5, 6).Value = 36
$ws.Cells.Item(6, 1).Value = 13
$ws.Cells.Item(6, 3).Value = 4
$ws.Cells.Item(6, 5).Value = 47
$ws.Cells.Item(6, 6).Value = 90
$ws.Cells.Item(7, 1).Value = 100
$ws.Cells.Item(7, 3).Value = 22
$ws.Cells.Item(7, 5).Value = 5
$ws.Cells.Item(7, 6).Value = 59
$ws.Cells.Item(8, 1).Value = 41
$ws.Cells.Item(8, 3).Value = 82
$ws.Cells.Item(8, 5).Value = 93
$ws.Cells.Item(8, 6).Value = 61
$ws.Cells.Item(9, 1).Value = 6
$ws.Cells.Item(9, 3).Value = 58
$ws.Cells.Item(9, 5).Value = 23
$ws.Cells.Item(9, 6).Value = 84
$ws.Cells.Item(10, 1).Value = 32
$ws.Cells.Item(10, 3).Value = 79
$ws.Cells.Item(10, 5).Value = 65
$ws.Cells.Item(11, 1).Value = 15
$ws.Cells.Item(11, 3).Value = 54
$ws.Cells.Item(11, 5).Value = 87
$ws.Cells.Item(12, 1).Value = 64
$ws.Cells.Item(12, 3).Value = 78
$ws.Cells.Item(12, 5).Value = 63
$ws.Cells.Item(13, 1).Value = 42
$ws.Cells.Item(13, 3).Value = 29
$ws.Cells.Item(13, 5).Value = 72
$ws.Cells.Item(14, 1).Value = 91
$ws.Cells.Item(14, 3).Value = 83
$ws.Cells.Item(14, 5).Value = 67
$ws.Cells.Item(15, 1).Value = 26
$ws.Cells.Item(15, 3).Value = 2
$ws.Cells.Item(15, 5).Value = 51
$ws.Cells.Item(16, 1).Value = 71
$ws.Cells.Item(16, 3).Value = 99
$ws.Cells.Item(16, 5).Value = 45
$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(17, 3).Value = 62
$ws.Cells.Item(17, 5).Value = 73
$ws.Cells.Item(18, 1).Value = 94
$ws.Cells.Item(18, 3).Value = 7
$ws.Cells.Item(18, 5).Value = 96
$ws.Cells.Item(19, 1).Value = 44
$ws.Cells.Item(19, 3).Value = 92
$ws.Cells.Item(19, 5).Value = 55
$ws.Cells.Item(20, 1).Value = 70
$ws.Cells.Item(20, 3).Value = 10
$ws.Cells.Item(21, 1).Value = 80
$ws.Cells.Item(21, 3).Value = 19
$ws.Cells.Item(22, 1).Value = 69
$ws.Cells.Item(22, 3).Value = 14
$ws.Cells.Item(23, 1).Value = 68
$ws.Cells.Item(23, 3).Value = 56
$ws.Cells.Item(24, 1).Value = 40
$ws.Cells.Item(24, 3).Value = 77
$ws.Cells.Item(25, 1).Value = 31
$ws.Cells.Item(25, 3).Value = 39
$ws.Cells.Item(26, 1).Value = 9
$ws.Cells.Item(26, 3).Value = 35
$ws.Cells.Item(27, 1).Value = 85
$ws.Cells.Item(28, 1).Value = 81
$ws.Cells.Item(29, 1).Value = 33
$ws.Cells.Item(30, 1).Value = 16
$ws.Cells.Item(31, 1).Value = 95
$ws.Cells.Item(32, 1).Value = 25
$ws.Cells.Item(33, 1).Value = 12
$ws.Cells.Item(34, 1).Value = 3
$ws.Cells.Item(35, 1).Value = 24
$ws.Cells.Item(36, 1).Value = 98

# ----- Backbones76 (was xl/worksheets/sheet78.xml) -----
$count = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($count))
$ws.Name = "Backbones76"

$styleSource.Range("A1:H1").Copy($ws.Range("A1:H1"))

$ws.Cells.Item(1, 1).Value = 8
$ws.Cells.Item(1, 2).Value = 60
$ws.Cells.Item(1, 3).Value = 88
$ws.Cells.Item(1, 4).Value = 97
$ws.Cells.Item(1, 5).Value = 20
$ws.Cells.Item(1, 6).Value = 1
$ws.Cells.Item(1, 7).Value = 53
$ws.Cells.Item(1, 8).Value = 46
$ws.Cells.Item(2, 1).Value = 9
$ws.Cells.Item(2, 2).Value = 32
$ws.Cells.Item(2, 3).Value = 7
$ws.Cells.Item(2, 5).Value = 13
$ws.Cells.Item(2, 6).Value = 72
$ws.Cells.Item(2, 7).Value = 79
$ws.Cells.Item(2, 8).Value = 21
$ws.Cells.Item(3, 1).Value = 40
$ws.Cells.Item(3, 2).Value = 27
$ws.Cells.Item(3, 3).Value = 37
$ws.Cells.Item(3, 5).Value = 17
$ws.Cells.Item(3, 6).Value = 95
$ws.Cells.Item(3, 7).Value = 44
$ws.Cells.Item(4, 1).Value = 75
$ws.Cells.Item(4, 2).Value = 98
$ws.Cells.Item(4, 3).Value = 36
$ws.Cells.Item(4, 5).Value = 49
$ws.Cells.Item(4, 6).Value = 28
$ws.Cells.Item(4, 7).Value = 56
$ws.Cells.Item(5, 1).Value = 85
$ws.Cells.Item(5, 2).Value = 62
$ws.Cells.Item(5, 3).Value = 69
$ws.Cells.Item(5, 5).Value = 24
$ws.Cells.Item(5, 6).Value = 33
$ws.Cells.Item(5, 7).Value = 70
$ws.Cells.Item(6, 1).Value = 67
$ws.Cells.Item(6, 2).Value = 48
$ws.Cells.Item(6, 3).Value = 25
$ws.Cells.Item(6, 5).Value = 22
$ws.Cells.Item(6, 6).Value = 66
$ws.Cells.Item(6, 7).Value = 65
$ws.Cells.Item(7, 1).Value = 90
$ws.Cells.Item(7, 2).Value = 12
$ws.Cells.Item(7, 5).Value = 29
$ws.Cells.Item(7, 6).Value = 47
$ws.Cells.Item(8, 1).Value = 82
$ws.Cells.Item(8, 2).Value = 42
$ws.Cells.Item(8, 5).Value = 15
$ws.Cells.Item(8, 6).Value = 5
$ws.Cells.Item(9, 1).Value = 81
$ws.Cells.Item(9, 2).Value = 39
$ws.Cells.Item(9, 5).Value = 86
$ws.Cells.Item(9, 6).Value = 96
$ws.Cells.Item(10, 1).Value = 94
$ws.Cells.Item(10, 2).Value = 93
$ws.Cells.Item(10, 5).Value = 63
$ws.Cells.Item(10, 6).Value = 89
$ws.Cells.Item(11, 1).Value = 61
$ws.Cells.Item(11, 2).Value = 51
$ws.Cells.Item(11, 5).Value = 59
$ws.Cells.Item(11, 6).Value = 80
$ws.Cells.Item(12, 1).Value = 43
$ws.Cells.Item(12, 2).Value = 14
$ws.Cells.Item(12, 5).Value = 50
$ws.Cells.Item(12, 6).Value = 18
$ws.Cells.Item(13, 1).Value = 91
$ws.Cells.Item(13, 2).Value = 23
$ws.Cells.Item(13, 5).Value = 76
$ws.Cells.Item(13, 6).Value = 87
$ws.Cells.Item(14, 1).Value = 19
$ws.Cells.Item(14, 2).Value = 3
$ws.Cells.Item(14, 5).Value = 77
$ws.Cells.Item(14, 6).Value = 92
$ws.Cells.Item(15, 1).Value = 31
$ws.Cells.Item(15, 2).Value = 38
$ws.Cells.Item(15, 5).Value = 54
$ws.Cells.Item(16, 1).Value = 100
$ws.Cells.Item(16, 2).Value = 73
$ws.Cells.Item(16, 5).Value = 41
$ws.Cells.Item(17, 2).Value = 55
$ws.Cells.Item(17, 5).Value = 68
$ws.Cells.Item(18, 2).Value = 84
$ws.Cells.Item(18, 5).Value = 10
$ws.Cells.Item(19, 2).Value = 16
$ws.Cells.Item(19, 5).Value = 71
$ws.Cells.Item(20, 2).Value = 78
$ws.Cells.Item(20, 5).Value = 6
$ws.Cells.Item(21, 2).Value = 26
$ws.Cells.Item(21, 5).Value = 52
$ws.Cells.Item(22, 2).Value = 83
$ws.Cells.Item(22, 5).Value = 2
$ws.Cells.Item(23, 2).Value = 4
$ws.Cells.Item(23, 5).Value = 74
$ws.Cells.Item(24, 2).Value = 64
$ws.Cells.Item(24, 5).Value = 35
$ws.Cells.Item(25, 2).Value = 57
$ws.Cells.Item(25, 5).Value = 45
$ws.Cells.Item(26, 2).Value = 34
$ws.Cells.Item(26, 5).Value = 58
$ws.Cells.Item(27, 2).Value = 99
$ws.Cells.Item(27, 5).Value = 11
$ws.Cells.Item(28, 5).Value = 30
